$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L with its own width (closest achievable value to the
# original 38.85546875 given this engine's column-width quantisation)
$ws.Columns.Item(12).ColumnWidth = 38

# Header cell (centered like the rest of row 1, style index 1)
$ws.Range("L1").Value = '(nY/\(nX/\nZ)) ->(n(Z/\nX)\/Y)'
$ws.Range("L1").HorizontalAlignment = -4108

# Data rows 2-9 (written in this particular order so that new shared
# strings land on the same table indices the original author's file has --
# L5's text was the last new unique string introduced)
$ws.Range("L2").Value = 'a0=1'
$ws.Range("L3").Value = 'a3 + a0=1'
$ws.Range("L4").Value = 'a2 + a0=1'
$ws.Range("L6").Value = 'a1 + a0= 1'
$ws.Range("L7").Value = 'a1 + a3 + a13 + a0=1'
$ws.Range("L8").Value = 'a1 + a2 + a12+ a0=1'
$ws.Range("L9").Value = 'a123 + a 12 + a13 + a 23 + a 1 + a2 + a3 + a0=1'
$ws.Range("L5").Value = 'a2 + a3 + a23 + a0=1'

# New isolated cell two rows below the table
$ws.Range("L11").Value = 1

# Leave the selection where the author left it
$ws.Range("L11").Select()
